$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - index 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 8217
$ws1.Range("G2").Value = 88
$ws1.Range("F3").Value = 128
$ws1.Range("F4").Value = 98
$ws1.Range("F5").Value = 33421
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 609
$ws1.Range("F8").Value = 716
$ws1.Range("F10").Value = 147
$ws1.Range("F11").Value = 450
$ws1.Range("F12").Value = 812
$ws1.Range("F13").Value = 62
$ws1.Range("F14").Value = 626
$ws1.Range("F15").Value = 429
$ws1.Range("F17").Value = 571
$ws1.Range("F20").Value = 425
$ws1.Range("F23").Value = 732
$ws1.Range("F24").Value = 2380
$ws1.Range("F25").Value = 858
$ws1.Range("F29").Value = 653
$ws1.Range("F30").Value = 653
$ws1.Range("F31").Value = 9
$ws1.Range("F32").Value = 1089

# Sheet "演出" (sheet2) - index 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 353

# Sheet "本地生活" (sheet3) - index 3
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 542

# Sheet "全部类型" (sheet4) - index 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 542
$ws4.Range("F3").Value = 8217
$ws4.Range("G3").Value = 88
$ws4.Range("F4").Value = 128
$ws4.Range("F5").Value = 98
$ws4.Range("F7").Value = 33421
$ws4.Range("F8").Value = 50
$ws4.Range("F9").Value = 609
$ws4.Range("F10").Value = 716
$ws4.Range("F13").Value = 147
$ws4.Range("F14").Value = 450
$ws4.Range("F15").Value = 353
$ws4.Range("F18").Value = 812
$ws4.Range("F19").Value = 62
$ws4.Range("F20").Value = 626
$ws4.Range("F21").Value = 429
$ws4.Range("F28").Value = 571
$ws4.Range("F31").Value = 425
$ws4.Range("F34").Value = 732
$ws4.Range("F35").Value = 2380
$ws4.Range("F36").Value = 858
$ws4.Range("F41").Value = 653
$ws4.Range("F42").Value = 653
$ws4.Range("F43").Value = 9
$ws4.Range("F44").Value = 1089

$wb.Save()
